# Update "Översikt HELSINGBORG" workbook to reflect the latest data refresh:
#  - bump the "Förändrad" (changed) date from 2023-09-23 (45192) to 2023-10-03 (45202)
#    for every existing record (rows 2-12)
#  - refresh the species/count totals on row 2 (a new inventory was added)
#  - add two newly observed species to row 2's species list
#  - append a brand-new logging notification record as row 13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: refreshed counts + species list ----
$ws.Range("C2").Value = 45202
$ws.Range("H2").Value = 14
$ws.Range("J2").Value = 21
$ws.Range("O2").Value = 26
$ws.Range("Q2").Value = 53

$speciesList = @(
  "Skogsalm","Ask","Blåsvart brunbagge","Gråtrut","Skillerticka","Bokvårtlav",
  "Dvärgklokrypare","Dvärgstumpbagge","Fyrfläckad vedsvampbagge","Getinglik svampmygga",
  "Gråbandad trägnagare","Gulsparv","Havsörn","Hålskenknäppare","Kråka","Matt pricklav",
  "Orangepudrad klotterlav","Plectophloeus nubigena","Skogsveronika","Skånebjörnbär",
  "Stiftklotterlav","Synchita variegata","Trinodes hirtus","Tvåfärgad barksvartbagge",
  "Vasstandad trädbasbagge","Violettgrå porlav","Blanksvart trämyra","Blåmossa",
  "Grå skärelav","Guldlockmossa","Gulnål","Havstulpanlav","Hässleklocka","Lundvårlök",
  "Noshornsoxe","Platt fjädermossa","Robust tickgnagare","Skogsbräsma","Skuggsprötmossa",
  "Stor häxört","Stubbspretmossa","Trubbfjädermossa","Tvåblad","Västlig hakmossa",
  "Dvärgpipistrell","Större brunfladdermus","Åkergroda","Vanlig groda","Vanlig padda",
  "Getlav","Grönvit nattviol","Sankt pers nycklar","Blåsippa"
)
$ws.Range("R2").Value = [string]::Join("`r`n", $speciesList)
$ws.Rows.Item(2).RowHeight = 15

# ---- Rows 3-12: same "Förändrad" date refresh ----
$ws.Range("C3").Value = 45202
$ws.Range("C4").Value = 45202
$ws.Range("C5").Value = 45202
$ws.Range("C6").Value = 45202
$ws.Range("C7").Value = 45202
$ws.Range("C8").Value = 45202
$ws.Range("C9").Value = 45202
$ws.Range("C10").Value = 45202
$ws.Range("C11").Value = 45202
$ws.Range("C12").Value = 45202

# Row 12 gets the same explicit row height as the rest of the data rows.
$ws.Rows.Item(12).RowHeight = 15

# ---- Row 13: brand-new logging notification ----
$ws.Range("A13").Value = "A 45832-2023"
$ws.Range("B13").Value = 45195
$ws.Range("C13").Value = 45202
$ws.Range("B13:C13").NumberFormat = $ws.Range("B12").NumberFormat
$ws.Range("D13").Value = "SKÅNE LÄN"
$ws.Range("E13").Value = "HELSINGBORG"
$ws.Range("G13").Value = 2.3
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0
$ws.Range("R13").Value = ""
$ws.Range("R13").WrapText = $true
